$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "2024-06-15 03:13:43"
$ws.Range("D12").Value = 200
$ws.Range("E12").Value = 2

# Row 13
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "2024-06-15 03:13:43"
$ws.Range("D13").Value = 200
$ws.Range("E13").Value = 0
